$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.231.75"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.905.03"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5354"
$ws.Range("E7").Value = "  +2.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3821"
$ws.Range("E8").Value = "  +1.55%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07300"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("E10").Value = "  +5.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9060"
$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08210"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.76"
$ws.Range("E13").Value = "  -0.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.353"
$ws.Range("E14").Value = "  +1.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.001"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.86"
$ws.Range("E16").Value = "  +2.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008671"
$ws.Range("E17").Value = "  +0.72%  "

$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").Value = "27.255.67"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.053"

$ws.Range("D21").Value = "1.072.19"
$ws.Range("E21").Value = "  -43.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.523"
$ws.Range("E23").Value = "  +1.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.09"
$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.295"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.38"
$ws.Range("E26").Value = "  +0.88%  "

$ws.Range("E27").Value = "  -0.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.83"
$ws.Range("E28").Value = "  +1.38%  "

$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.726"
$ws.Range("E30").Value = "  -4.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09222"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8302"
$ws.Range("E32").Value = "  +4.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05082"

$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.006"
$ws.Range("E35").Value = "  +1.86%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.323"
$ws.Range("E36").Value = "  -3.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.677"
$ws.Range("E37").Value = "  +3.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5903"
$ws.Range("E38").Value = "  +4.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02003"
$ws.Range("E39").Value = "  +0.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.078"
$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.355"
$ws.Range("E41").Value = "  +4.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.640"
$ws.Range("E42").Value = "  +1.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.30"
$ws.Range("E43").Value = "  +1.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5099"
$ws.Range("E44").Value = "  +4.08%  "

$ws.Range("E45").Value = "  +0.63%  "

$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.17"
$ws.Range("E47").Value = "  +1.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.37"
$ws.Range("E49").Value = "  +0.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06159"
$ws.Range("E50").Value = "  +3.62%  "

$ws.Range("E51").Value = "  +0.15%  "
